$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.676.54'
$ws.Range('E2').Value = '  -0.10%  '

$ws.Range('D3').Value = '3.527.19'
$ws.Range('E3').Value = '  -1.21%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '625.00'
$ws.Range('E5').Value = '  +2.95%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.62'
$ws.Range('E6').Value = '  -0.34%  '

$ws.Range('D7').Value = '3.527.70'
$ws.Range('E7').Value = '  -1.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.610'
$ws.Range('E8').Value = '  -0.95%  '

$ws.Range('E9').Value = '  +0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.199'
$ws.Range('E10').Value = '  +1.28%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.15'
$ws.Range('E11').Value = '  -3.10%  '

$ws.Range('E12').Value = '  -0.55%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.42'
$ws.Range('E13').Value = '  -1.00%  '

$ws.Range('E14').Value = '  -0.44%  '

$ws.Range('D15').Value = '4.095.66'
$ws.Range('E15').Value = '  -1.15%  '

$ws.Range('E16').Value = '  -0.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '606.43'
$ws.Range('E17').Value = '  -1.76%  '

$ws.Range('D18').Value = '3.533.48'
$ws.Range('E18').Value = '  -0.93%  '

$ws.Range('D19').Value = '70.759.69'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.121'
$ws.Range('E20').Value = '  +1.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.75'
$ws.Range('E21').Value = '  +1.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.884'
$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.10'
$ws.Range('E23').Value = '  -3.06%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '98.16'
$ws.Range('E24').Value = '  +0.65%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.64'
$ws.Range('E25').Value = '  -2.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.78'
$ws.Range('E26').Value = '  -0.99%  '

$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.57'
$ws.Range('E28').Value = '  -2.68%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.74'
$ws.Range('E29').Value = '  +0.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.07'
$ws.Range('E30').Value = '  -0.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.03'
$ws.Range('E31').Value = '  -0.56%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.12'
$ws.Range('E32').Value = '  -4.28%  '

$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '634.80'
$ws.Range('E34').Value = '  +3.61%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.79'
$ws.Range('E35').Value = '  -3.39%  '

$ws.Range('B36').Value = 'Cosmos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.82'
$ws.Range('E36').Value = '  -0.34%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0997'
$ws.Range('E37').Value = '  -2.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.46'
$ws.Range('E38').Value = '  -7.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0473'
$ws.Range('E39').Value = '  -1.84%  '

$ws.Range('E40').Value = '  -1.00%  '

$ws.Range('E41').Value = '  +0.41%  '

$ws.Range('E42').Value = '  +1.46%  '

$ws.Range('D43').Value = '3.357.48'
$ws.Range('E43').Value = '  -0.92%  '

$ws.Range('D44').Value = '0.0₃0732'
$ws.Range('E44').Value = '  +3.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.98'
$ws.Range('E45').Value = '  -0.58%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.311'
$ws.Range('E46').Value = '  -3.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '32.02'
$ws.Range('E47').Value = '  -3.09%  '

$ws.Range('E48').Value = '  -2.51%  '

$ws.Range('E49').Value = '  -0.26%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.86'
$ws.Range('E50').Value = '  -0.11%  '

$ws.Range('E51').Value = '  +5.70%  '
